$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue {
    param($ws, $ref, $val)
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws 'D2' '318.41'
Set-TextValue $ws 'E2' '3.19%'
Set-TextValue $ws 'D3' '41.33'
Set-TextValue $ws 'E3' '0.70%'
Set-TextValue $ws 'D4' '5.237'
Set-TextValue $ws 'E4' '2.04%'
Set-TextValue $ws 'D5' '0.07727'
Set-TextValue $ws 'E5' '1.21%'
Set-TextValue $ws 'D6' '1.703'
Set-TextValue $ws 'E6' '4.94%'
Set-TextValue $ws 'D7' '0.9507'
Set-TextValue $ws 'E7' '4.54%'
Set-TextValue $ws 'E8' '-1.89%'
Set-TextValue $ws 'D9' '0.1266'
Set-TextValue $ws 'E9' '8.20%'
Set-TextValue $ws 'D10' '0.1825'
Set-TextValue $ws 'E10' '1.48%'
Set-TextValue $ws 'D11' '0.09206'
Set-TextValue $ws 'E11' '1.01%'
Set-TextValue $ws 'D12' '0.04366'
Set-TextValue $ws 'E12' '2.66%'
Set-TextValue $ws 'E13' '0.88%'
Set-TextValue $ws 'D14' '0.001259'
Set-TextValue $ws 'E14' '0.05%'
Set-TextValue $ws 'D15' '0.005871'
Set-TextValue $ws 'E15' '0.26%'
Set-TextValue $ws 'E16' '0.03%'
Set-TextValue $ws 'D17' '4.298'
Set-TextValue $ws 'E17' '0.44%'
Set-TextValue $ws 'E18' '2.88%'
Set-TextValue $ws 'D19' '7.494'
Set-TextValue $ws 'E19' '8.18%'
Set-TextValue $ws 'D20' '0.1348'
Set-TextValue $ws 'E20' '-3.09%'
Set-TextValue $ws 'D21' '0.2814'
Set-TextValue $ws 'E21' '4.00%'
Set-TextValue $ws 'D22' '0.04017'
Set-TextValue $ws 'E22' '-0.32%'
Set-TextValue $ws 'D23' '0.001263'
Set-TextValue $ws 'E23' '-0.68%'
Set-TextValue $ws 'D24' '0.004248'
Set-TextValue $ws 'E24' '4.12%'
Set-TextValue $ws 'D25' '0.0001270'
Set-TextValue $ws 'E25' '-0.18%'
Set-TextValue $ws 'D38' '0.02549'
Set-TextValue $ws 'E38' '5.27%'
Set-TextValue $ws 'D39' '0.05361'
Set-TextValue $ws 'E39' '2.54%'
Set-TextValue $ws 'D40' '0.007785'
Set-TextValue $ws 'E40' '-0.22%'
Set-TextValue $ws 'D41' '0.1321'
Set-TextValue $ws 'E41' '1.56%'
Set-TextValue $ws 'D42' '0.007327'
Set-TextValue $ws 'E42' '7.72%'
Set-TextValue $ws 'D43' '0.001975'
Set-TextValue $ws 'E43' '3.74%'
Set-TextValue $ws 'D44' '0.007605'
Set-TextValue $ws 'E44' '-5.80%'
Set-TextValue $ws 'D45' '0.3436'
Set-TextValue $ws 'E45' '12.03%'
Set-TextValue $ws 'D46' '0.00006710'
Set-TextValue $ws 'E46' '-2.73%'
Set-TextValue $ws 'D47' '0.00000000750'
Set-TextValue $ws 'E47' '-0.19%'
Set-TextValue $ws 'D48' '0.2229'
Set-TextValue $ws 'E48' '135.08%'
Set-TextValue $ws 'D49' '0.004188'
Set-TextValue $ws 'E49' '39.47%'
Set-TextValue $ws 'D50' '0.00002100'
Set-TextValue $ws 'E50' '-0.19%'
Set-TextValue $ws 'D51' '0.0002000'
Set-TextValue $ws 'E51' '-0.19%'
